# Add data for 2021-11-19
# Updates the "carjacking arrests by month" workbook so that the "as-of"
# date moves from 2021-11-10 to 2021-11-11, adding one more day's worth
# of arrest data to November 2021 and to the running totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet (and, implicitly, update the workbook tab title).
$ws.Name = "Through 2021-11-11"

# 2. Row 6 (April) - 2021 columns (T/U/V)
$ws.Range("T6").Value = 13
$ws.Range("U6").Value = 87
$ws.Range("V6").Value = 0.13

# 3. Row 12 (October) - 2021 columns (T/U/V)
$ws.Range("T12").Value = 4
$ws.Range("U12").Value = 192
$ws.Range("V12").Value = 0.0204

# 4. Row 13 (November) - label + most columns updated, plus brand new
#    T13/V13 cells (2021 arrest_made / arrest_rate for November).
$ws.Range("A13").Value = "November (through 11-11)"

$ws.Range("C13").Value = 13

$ws.Range("F13").Value = 25
$ws.Range("G13").Value = 0.0385

$ws.Range("I13").Value = 41
$ws.Range("J13").Value = 0.0238

$ws.Range("N13").Value = 3
$ws.Range("O13").Value = 16
$ws.Range("P13").Value = 0.1579

$ws.Range("R13").Value = 68
$ws.Range("S13").Value = 0.0286

$ws.Range("T13").Value = 1
$ws.Range("U13").Value = 73
$ws.Range("V13").Value = 0.0135
$ws.Range("V13").NumberFormat = $ws.Range("V12").NumberFormat

# 5. Row 14 (Total) - running totals shift for every year that had new
#    November activity.
$ws.Range("C14").Value = 239
$ws.Range("D14").Value = 0.1181

$ws.Range("F14").Value = 459
$ws.Range("G14").Value = 0.1035

$ws.Range("I14").Value = 690
$ws.Range("J14").Value = 0.0824

$ws.Range("N14").Value = 51
$ws.Range("O14").Value = 450
$ws.Range("P14").Value = 0.1018

$ws.Range("R14").Value = 1071
$ws.Range("S14").Value = 0.0497

$ws.Range("T14").Value = 89
$ws.Range("U14").Value = 1429
$ws.Range("V14").Value = 0.0586
